$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.610.82"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.811.72"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.91"
$ws.Range("E5").Value = "  -1.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.600"
$ws.Range("E6").Value = "  +3.23%  "

$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.57"
$ws.Range("E8").Value = "  +5.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.293"
$ws.Range("E9").Value = "  -2.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0684"
$ws.Range("E10").Value = "  -1.60%  "

$ws.Range("E11").Value = "  +1.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.072.47"
$ws.Range("E12").Value = "  +0.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.33"
$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.827.33"
$ws.Range("E14").Value = "  +1.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.632"
$ws.Range("E15").Value = "  -1.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.533.85"
$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.44"
$ws.Range("E17").Value = "  +1.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.49"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.33"
$ws.Range("E19").Value = "  -0.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = [string]::Concat("0.0", [char]0x2083, "0778")
$ws.Range("E20").Value = "  -2.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.25"
$ws.Range("E21").Value = "  -2.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.12"
$ws.Range("E23").Value = "  -1.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("E24").Value = "  +5.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.63"
$ws.Range("E25").Value = "  -1.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.87"
$ws.Range("E26").Value = "  +0.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.32"
$ws.Range("E27").Value = "  +3.17%  "

$ws.Range("E28").Value = "  +1.90%  "

$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.83"
$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.92"
$ws.Range("E32").Value = "  -2.06%  "

$ws.Range("E33").Value = "  -2.75%  "

$ws.Range("E34").Value = "  -1.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.363.88"
$ws.Range("E35").Value = "  -2.18%  "

$ws.Range("E36").Value = "  -4.13%  "

$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.37"
$ws.Range("E38").Value = "  -5.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0187"
$ws.Range("E39").Value = "  -1.74%  "

$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "81.35"
$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("E42").Value = "  -1.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.938"
$ws.Range("E43").Value = "  -1.18%  "

$ws.Range("E44").Value = "  +4.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.61"
$ws.Range("E45").Value = "  +0.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0501"
$ws.Range("E46").Value = "  -1.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.972.87"
$ws.Range("E47").Value = "  +0.46%  "

$ws.Range("E48").Value = "  -2.62%  "


$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.72"
$ws.Range("E50").Value = "  -2.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = [string]::Concat("0.0", [char]0x2086, "0121")
$ws.Range("E51").Value = "  -6.76%  "
